# Update odds data in Sheet1 for the 2024-11-08 FlashScore export.
# This reflects refreshed odds values scraped after the initial commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Internacional vs Fluminense ---
$ws.Range("G2").Value = 1.48
$ws.Range("H2").Value = 4.2
$ws.Range("J2").Value = 2.05
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 3.4
$ws.Range("Q2").Value = 2.05
$ws.Range("R2").Value = 1.75
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 9.5
$ws.Range("AG2").Value = 15
$ws.Range("BB2").Value = 451

# --- Row 3: Racing Montevideo vs Nacional ---
$ws.Range("G3").Value = 7
$ws.Range("H3").Value = 3.8
$ws.Range("I3").Value = 1.57
$ws.Range("L3").Value = 2.2
$ws.Range("AN3").Value = 8
$ws.Range("BA3").Value = 51
